$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The B column (month/year) values are being re-dated from 2024 to 2020,
# keeping the same month-of-year, and are switched to an explicit Text
# number format so Excel does not reinterpret "Mon/Year" as a date.
$months = @("Jan/2020","Feb/2020","Mar/2020","Apr/2020","May/2020","Jun/2020","Jul/2020","Aug/2020","Sep/2020","Oct/2020","Nov/2020","Dec/2020")

$range = $ws.Range("B2:B61")
$range.NumberFormat = "@"

for ($r = 2; $r -le 61; $r++) {
    $m = $months[($r - 2) % 12]
    $ws.Cells.Item($r, 2).Value = $m
}

# Best-effort cosmetic match: column widths to fit the refreshed content,
# and the active selection left on C4 as in the edited workbook.
$ws.Columns.Item(1).ColumnWidth = 6.85546875
$ws.Columns.Item(2).ColumnWidth = 9.5703125
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 17.5703125
$ws.Columns.Item(5).ColumnWidth = 12.85546875
$ws.Columns.Item(6).ColumnWidth = 18.42578125
$ws.Columns.Item(7).ColumnWidth = 16.140625
$ws.Columns.Item(8).ColumnWidth = 21.7109375

$ws.Range("C4").Select()
